$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D and E on these rows are stored as text (inlineStr) in the
# original workbook (e.g. "290.20", "-3.62%"), so force the Text number
# format before assigning to keep Excel from re-interpreting the values
# as numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "290.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.62%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.33%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.949"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.84%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07132"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-9.26%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.812"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-13.66%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.677"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.87%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.735"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.54%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8983"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.10%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1642"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.50%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07556"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.10%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07994"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.88%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03053"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.28%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.10%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001491"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.94%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005659"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.78%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.478"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.01%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.101"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.77%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3276"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.01%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1278"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.68%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.044"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.75%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.2050"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "10.49%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04518"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.51%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001212"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.93%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003991"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.26%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001249"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01610"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-8.41%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04360"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-8.99%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007347"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.56%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1304"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.12%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.52%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009248"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.93%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006032"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.17%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.247"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "173.87%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003000"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-11.43%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.35%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.35%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.35%"
